# "Add PPL data and targets" - one more row of data (PPL Corp.) added to
# the ITR input data, ITR target input data and Portfolio sheets, plus the
# selection/active-tab view state that Excel recorded when the author saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "ITR input data" (index 4) - new row 33 for PPL Corp.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Pull formatting for the brand-new cells (B, G, I, O, P, AE:AI) from the
# row above, which already carries the correct per-column styles.
$ws4.Range("A32:U32").Copy()
$ws4.Range("A33:U33").PasteSpecial(-4122)
$ws4.Range("AE32:AI32").Copy()
$ws4.Range("AE33:AI33").PasteSpecial(-4122)
$ws4.Range("AS32:AW32").Copy()
$ws4.Range("AS33:AW33").PasteSpecial(-4122)

$ws4.Cells.Item(33,1).Value  = "PPL Corp."
$ws4.Cells.Item(33,2).Value  = "9N3UAJSNOUXFKQLF3V18"
$ws4.Cells.Item(33,3).Value  = "US69351T1060"
$ws4.Cells.Item(33,4).Value  = "US"
$ws4.Cells.Item(33,5).Value  = "North America"
$ws4.Cells.Item(33,6).Value  = "Electricity Utilities"
$ws4.Cells.Item(33,7).Value  = "equity"
$ws4.Cells.Item(33,8).Value  = "USD"
$ws4.Cells.Item(33,9).Value  = 44196
$ws4.Cells.Item(33,10).Value = 19865342074
$ws4.Cells.Item(33,11).Value = 7769000000
$ws4.Cells.Item(33,12).Value = 40943342074
$ws4.Cells.Item(33,13).Value = 41758342074
$ws4.Cells.Item(33,14).Value = 45680000000
$ws4.Cells.Item(33,15).Value = "Mt CO2"
$ws4.Cells.Item(33,16).Value = "TWh"

$ws4.Cells.Item(33,31).Value = 30.088487220000001
$ws4.Cells.Item(33,32).Value = 30.24837145
$ws4.Cells.Item(33,33).Value = 31.611469039999999
$ws4.Cells.Item(33,34).Value = 28.778915319999999
$ws4.Cells.Item(33,35).Value = 28.07780713

$ws4.Cells.Item(33,45).Value = 38.355258640000002
$ws4.Cells.Item(33,46).Value = 37.442832350000003
$ws4.Cells.Item(33,47).Value = 39.590075179999999
$ws4.Cells.Item(33,48).Value = 35.152931719999998
$ws4.Cells.Item(33,49).Value = 32.487984334642732

# ---------------------------------------------------------------------
# Sheet "ITR target input data" (index 5) - two new target rows for PPL
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("A41:H41").Copy()
$ws5.Range("A42:H42").PasteSpecial(-4122)
$ws5.Range("J41:L41").Copy()
$ws5.Range("J42:L42").PasteSpecial(-4122)

$ws5.Cells.Item(42,1).Value = "PPL Corp."
$ws5.Cells.Item(42,2).Value = "9N3UAJSNOUXFKQLF3V18"
$ws5.Cells.Item(42,3).Value = "US69351T1060"
$ws5.Cells.Item(42,4).Value = 2050
$ws5.Cells.Item(42,5).Value = "absolute"
$ws5.Cells.Item(42,6).Value = "S1+S2"
$ws5.Cells.Item(42,7).Value = 2021
$ws5.Cells.Item(42,8).Value = 2010
$ws5.Cells.Item(42,9).Formula = "=60736086+1597157"
$ws5.Cells.Item(42,10).Value = "t CO2"
$ws5.Cells.Item(42,11).Value = 2035
$ws5.Cells.Item(42,12).Value = 0.7

$ws5.Range("A41:H41").Copy()
$ws5.Range("A43:H43").PasteSpecial(-4122)
$ws5.Range("J41:L41").Copy()
$ws5.Range("J43:L43").PasteSpecial(-4122)

$ws5.Cells.Item(43,1).Value = "PPL Corp."
$ws5.Cells.Item(43,2).Value = "9N3UAJSNOUXFKQLF3V18"
$ws5.Cells.Item(43,3).Value = "US69351T1060"
$ws5.Cells.Item(43,4).Value = 2050
$ws5.Cells.Item(43,5).Value = "absolute"
$ws5.Cells.Item(43,6).Value = "S1+S2"
$ws5.Cells.Item(43,7).Value = 2021
$ws5.Cells.Item(43,8).Value = 2010
$ws5.Cells.Item(43,9).Formula = "=60736086+1597157"
$ws5.Cells.Item(43,10).Value = "t CO2"
$ws5.Cells.Item(43,11).Value = 2040
$ws5.Cells.Item(43,12).Value = 0.8

# ---------------------------------------------------------------------
# Sheet "Portfolio" (index 7) - new row 33 for PPL Corp.
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)

$ws7.Range("A3:D3").Copy()
$ws7.Range("A33:D33").PasteSpecial(-4122)

$ws7.Cells.Item(33,1).Value = "PPL Corp."
$ws7.Cells.Item(33,2).Value = "9N3UAJSNOUXFKQLF3V18"
$ws7.Cells.Item(33,3).Value = "US69351T1060"
$ws7.Cells.Item(33,4).Value = "US69351T1060"
$ws7.Cells.Item(33,5).Formula = "=RANDBETWEEN(35000,250000)"

# ---------------------------------------------------------------------
# View state: selections on each touched sheet, and which sheet/cell was
# left selected when the workbook was last saved (Portfolio / E33).
# ---------------------------------------------------------------------
$ws4.Activate()
$ws4.Range("A33:C33").Select()

$ws5.Activate()
$ws5.Range("A43").Select()

$ws8 = $wb.Worksheets.Item(8)
$ws8.Activate()
$ws8.Rows.Item(40).Select()

$ws9 = $wb.Worksheets.Item(9)
$ws9.Activate()
$ws9.Range("F4").Select()

$ws7.Activate()
$ws7.Range("E33").Select()
